$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.578493666666667
$ws.Range("H2").Value = 16.735481
$ws.Range("I2").Value = 0.1036332930693284
$ws.Range("J2").Value = 0.1036332930693284
$ws.Range("M2").Value = 2.324047
$ws.Range("N2").Value = 6.972140999999999
$ws.Range("O2").Value = 0.2694410417186929
$ws.Range("P2").Value = 0.2694410417186929
$ws.Range("Q2").Value = 12.96468147053567
$ws.Range("R2").Value = 116.682133234821
$ws.Range("S2").Value = 0.02792306244133844
$ws.Range("T2").Value = 0.02792306244133845

$ws.Range("G3").Value = 5.578493666666667
$ws.Range("H3").Value = 16.735481
$ws.Range("I3").Value = 0.1036332930693284
$ws.Range("J3").Value = 0.1036332930693284
$ws.Range("O3").Value = 0.1652951421133628
$ws.Range("P3").Value = 0.1652951421133628
$ws.Range("Q3").Value = 7.953498295794334
$ws.Range("R3").Value = 71.58148466214901
$ws.Range("S3").Value = 0.01713007990557041
$ws.Range("T3").Value = 0.01713007990557042

$ws.Range("G4").Value = 5.578493666666667
$ws.Range("H4").Value = 16.735481
$ws.Range("I4").Value = 0.1036332930693284
$ws.Range("J4").Value = 0.1036332930693284
$ws.Range("M4").Value = 2.77793
$ws.Range("N4").Value = 8.33379
$ws.Range("O4").Value = 0.3220624854065381
$ws.Range("P4").Value = 0.3220624854065382
$ws.Range("Q4").Value = 15.49666491144333
$ws.Range("R4").Value = 139.46998420299
$ws.Range("S4").Value = 0.03337639593677207
$ws.Range("T4").Value = 0.03337639593677208

$ws.Range("G5").Value = 5.578493666666667
$ws.Range("H5").Value = 16.735481
$ws.Range("I5").Value = 0.1036332930693284
$ws.Range("J5").Value = 0.1036332930693284
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4528016666666666
$ws.Range("N5").Value = 1.358405
$ws.Range("O5").Value = 0.05249607807356178
$ws.Range("P5").Value = 0.05249607807356178
$ws.Range("Q5").Value = 2.525951229756111
$ws.Range("R5").Value = 22.733561067805
$ws.Range("S5").Value = 0.005440341443987771
$ws.Range("T5").Value = 0.005440341443987774

$ws.Range("G6").Value = 5.578493666666667
$ws.Range("H6").Value = 16.735481
$ws.Range("I6").Value = 0.1036332930693284
$ws.Range("J6").Value = 0.1036332930693284
$ws.Range("M6").Value = 1.644916333333333
$ws.Range("N6").Value = 4.934749
$ws.Range("O6").Value = 0.1907052526878442
$ws.Range("P6").Value = 0.1907052526878442
$ws.Range("Q6").Value = 9.176155347696556
$ws.Range("R6").Value = 82.585398129269
$ws.Range("S6").Value = 0.01976341334165968
$ws.Range("T6").Value = 0.01976341334165969

$ws.Range("I7").Value = 0.06881911773528272
$ws.Range("J7").Value = 0.06881911773528274
$ws.Range("M7").Value = 2.324047
$ws.Range("N7").Value = 6.972140999999999
$ws.Range("O7").Value = 0.2694410417186929
$ws.Range("P7").Value = 0.2694410417186929
$ws.Range("Q7").Value = 8.609375559689665
$ws.Range("R7").Value = 77.48438003720699
$ws.Range("S7").Value = 0.01854269477275595
$ws.Range("T7").Value = 0.01854269477275595

$ws.Range("I8").Value = 0.06881911773528272
$ws.Range("J8").Value = 0.06881911773528274
$ws.Range("O8").Value = 0.1652951421133628
$ws.Range("P8").Value = 0.1652951421133628
$ws.Range("S8").Value = 0.0113754658461698
$ws.Range("T8").Value = 0.01137546584616981

$ws.Range("I9").Value = 0.06881911773528272
$ws.Range("J9").Value = 0.06881911773528274
$ws.Range("M9").Value = 2.77793
$ws.Range("N9").Value = 8.33379
$ws.Range("O9").Value = 0.3220624854065381
$ws.Range("P9").Value = 0.3220624854065382
$ws.Range("Q9").Value = 10.29077408870333
$ws.Range("R9").Value = 92.61696679833
$ws.Range("S9").Value = 0.02216405610131032
$ws.Range("T9").Value = 0.02216405610131033

$ws.Range("I10").Value = 0.06881911773528272
$ws.Range("J10").Value = 0.06881911773528274
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4528016666666666
$ws.Range("N10").Value = 1.358405
$ws.Range("O10").Value = 0.05249607807356178
$ws.Range("P10").Value = 0.05249607807356178
$ws.Range("Q10").Value = 1.677392755992777
$ws.Range("R10").Value = 15.096534803935
$ws.Range("S10").Value = 0.003612733777585042
$ws.Range("T10").Value = 0.003612733777585043

$ws.Range("I11").Value = 0.06881911773528272
$ws.Range("J11").Value = 0.06881911773528274
$ws.Range("M11").Value = 1.644916333333333
$ws.Range("N11").Value = 4.934749
$ws.Range("O11").Value = 0.1907052526878442
$ws.Range("P11").Value = 0.1907052526878442
$ws.Range("Q11").Value = 6.09355253053589
$ws.Range("R11").Value = 54.841972774823
$ws.Range("S11").Value = 0.01312416723746159
$ws.Range("T11").Value = 0.01312416723746159

$ws.Range("G12").Value = 24.77295966666667
$ws.Range("H12").Value = 74.31887900000001
$ws.Range("I12").Value = 0.4602144490493554
$ws.Range("J12").Value = 0.4602144490493556
$ws.Range("M12").Value = 2.324047
$ws.Range("N12").Value = 6.972140999999999
$ws.Range("O12").Value = 0.2694410417186929
$ws.Range("P12").Value = 0.2694410417186929
$ws.Range("Q12").Value = 57.57352259443766
$ws.Range("R12").Value = 518.161703349939
$ws.Range("S12").Value = 0.1240006605658527
$ws.Range("T12").Value = 0.1240006605658527

$ws.Range("G13").Value = 24.77295966666667
$ws.Range("H13").Value = 74.31887900000001
$ws.Range("I13").Value = 0.4602144490493554
$ws.Range("J13").Value = 0.4602144490493556
$ws.Range("O13").Value = 0.1652951421133628
$ws.Range("P13").Value = 0.1652951421133628
$ws.Range("Q13").Value = 35.31987383403234
$ws.Range("R13").Value = 317.8788645062911
$ws.Range("S13").Value = 0.07607121275823617
$ws.Range("T13").Value = 0.07607121275823621

$ws.Range("G14").Value = 24.77295966666667
$ws.Range("H14").Value = 74.31887900000001
$ws.Range("I14").Value = 0.4602144490493554
$ws.Range("J14").Value = 0.4602144490493556
$ws.Range("M14").Value = 2.77793
$ws.Range("N14").Value = 8.33379
$ws.Range("O14").Value = 0.3220624854065381
$ws.Range("P14").Value = 0.3220624854065382
$ws.Range("Q14").Value = 68.81754784682335
$ws.Range("R14").Value = 619.3579306214101
$ws.Range("S14").Value = 0.148217809280836
$ws.Range("T14").Value = 0.1482178092808361

$ws.Range("G15").Value = 24.77295966666667
$ws.Range("H15").Value = 74.31887900000001
$ws.Range("I15").Value = 0.4602144490493554
$ws.Range("J15").Value = 0.4602144490493556
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.4528016666666666
$ws.Range("N15").Value = 1.358405
$ws.Range("O15").Value = 0.05249607807356178
$ws.Range("P15").Value = 0.05249607807356178
$ws.Range("Q15").Value = 11.21723742533278
$ws.Range("R15").Value = 100.955136827995
$ws.Range("S15").Value = 0.02415945364787618
$ws.Range("T15").Value = 0.02415945364787619

$ws.Range("G16").Value = 24.77295966666667
$ws.Range("H16").Value = 74.31887900000001
$ws.Range("I16").Value = 0.4602144490493554
$ws.Range("J16").Value = 0.4602144490493556
$ws.Range("M16").Value = 1.644916333333333
$ws.Range("N16").Value = 4.934749
$ws.Range("O16").Value = 0.1907052526878442
$ws.Range("P16").Value = 0.1907052526878442
$ws.Range("Q16").Value = 40.7494459807079
$ws.Range("R16").Value = 366.7450138263711
$ws.Range("S16").Value = 0.08776531279655431
$ws.Range("T16").Value = 0.08776531279655435

$ws.Range("G17").Value = 0.4291063333333334
$ws.Range("H17").Value = 1.287319
$ws.Range("I17").Value = 0.007971632676749163
$ws.Range("J17").Value = 0.007971632676749165
$ws.Range("M17").Value = 2.324047
$ws.Range("N17").Value = 6.972140999999999
$ws.Range("O17").Value = 0.2694410417186929
$ws.Range("P17").Value = 0.2694410417186929
$ws.Range("Q17").Value = 0.9972632866643333
$ws.Range("R17").Value = 8.975369579978999
$ws.Range("S17").Value = 0.002147885012622067
$ws.Range("T17").Value = 0.002147885012622067

$ws.Range("G18").Value = 0.4291063333333334
$ws.Range("H18").Value = 1.287319
$ws.Range("I18").Value = 0.007971632676749163
$ws.Range("J18").Value = 0.007971632676749165
$ws.Range("O18").Value = 0.1652951421133628
$ws.Range("P18").Value = 0.1652951421133628
$ws.Range("Q18").Value = 0.6117953510056667
$ws.Range("R18").Value = 5.506158159051001
$ws.Range("S18").Value = 0.00131767215617878
$ws.Range("T18").Value = 0.00131767215617878

$ws.Range("G19").Value = 0.4291063333333334
$ws.Range("H19").Value = 1.287319
$ws.Range("I19").Value = 0.007971632676749163
$ws.Range("J19").Value = 0.007971632676749165
$ws.Range("M19").Value = 2.77793
$ws.Range("N19").Value = 8.33379
$ws.Range("O19").Value = 0.3220624854065381
$ws.Range("P19").Value = 0.3220624854065382
$ws.Range("Q19").Value = 1.192027356556667
$ws.Range("R19").Value = 10.72824620901
$ws.Range("S19").Value = 0.00256736383262181
$ws.Range("T19").Value = 0.002567363832621811

$ws.Range("G20").Value = 0.4291063333333334
$ws.Range("H20").Value = 1.287319
$ws.Range("I20").Value = 0.007971632676749163
$ws.Range("J20").Value = 0.007971632676749165
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.4528016666666666
$ws.Range("N20").Value = 1.358405
$ws.Range("O20").Value = 0.05249607807356178
$ws.Range("P20").Value = 0.05249607807356178
$ws.Range("Q20").Value = 0.1943000629105555
$ws.Range("R20").Value = 1.748700566195
$ws.Range("S20").Value = 0.0004184794513723803
$ws.Range("T20").Value = 0.0004184794513723805

$ws.Range("G21").Value = 0.4291063333333334
$ws.Range("H21").Value = 1.287319
$ws.Range("I21").Value = 0.007971632676749163
$ws.Range("J21").Value = 0.007971632676749165
$ws.Range("M21").Value = 1.644916333333333
$ws.Range("N21").Value = 4.934749
$ws.Range("O21").Value = 0.1907052526878442
$ws.Range("P21").Value = 0.1907052526878442
$ws.Range("Q21").Value = 0.7058440164367779
$ws.Range("R21").Value = 6.352596147931001
$ws.Range("S21").Value = 0.001520232223954125
$ws.Range("T21").Value = 0.001520232223954125

$ws.Range("G22").Value = 19.34413
$ws.Range("H22").Value = 58.03239
$ws.Range("I22").Value = 0.3593615074692841
$ws.Range("J22").Value = 0.3593615074692842
$ws.Range("M22").Value = 2.324047
$ws.Range("N22").Value = 6.972140999999999
$ws.Range("O22").Value = 0.2694410417186929
$ws.Range("P22").Value = 0.2694410417186929
$ws.Range("Q22").Value = 44.95666729411
$ws.Range("R22").Value = 404.6100056469899
$ws.Range("S22").Value = 0.09682673892612377
$ws.Range("T22").Value = 0.09682673892612378

$ws.Range("G23").Value = 19.34413
$ws.Range("H23").Value = 58.03239
$ws.Range("I23").Value = 0.3593615074692841
$ws.Range("J23").Value = 0.3593615074692842
$ws.Range("O23").Value = 0.1652951421133628
$ws.Range("P23").Value = 0.1652951421133628
$ws.Range("Q23").Value = 27.57975793859
$ws.Range("R23").Value = 248.21782144731
$ws.Range("S23").Value = 0.05940071144720761
$ws.Range("T23").Value = 0.05940071144720763

$ws.Range("G24").Value = 19.34413
$ws.Range("H24").Value = 58.03239
$ws.Range("I24").Value = 0.3593615074692841
$ws.Range("J24").Value = 0.3593615074692842
$ws.Range("M24").Value = 2.77793
$ws.Range("N24").Value = 8.33379
$ws.Range("O24").Value = 0.3220624854065381
$ws.Range("P24").Value = 0.3220624854065382
$ws.Range("Q24").Value = 53.7366390509
$ws.Range("R24").Value = 483.6297514581
$ws.Range("S24").Value = 0.1157368602549979
$ws.Range("T24").Value = 0.1157368602549979

$ws.Range("G25").Value = 19.34413
$ws.Range("H25").Value = 58.03239
$ws.Range("I25").Value = 0.3593615074692841
$ws.Range("J25").Value = 0.3593615074692842
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.4528016666666666
$ws.Range("N25").Value = 1.358405
$ws.Range("O25").Value = 0.05249607807356178
$ws.Range("P25").Value = 0.05249607807356178
$ws.Range("Q25").Value = 8.759054304216665
$ws.Range("R25").Value = 78.83148873795
$ws.Range("S25").Value = 0.01886506975274039
$ws.Range("T25").Value = 0.0188650697527404

$ws.Range("G26").Value = 19.34413
$ws.Range("H26").Value = 58.03239
$ws.Range("I26").Value = 0.3593615074692841
$ws.Range("J26").Value = 0.3593615074692842
$ws.Range("M26").Value = 1.644916333333333
$ws.Range("N26").Value = 4.934749
$ws.Range("O26").Value = 0.1907052526878442
$ws.Range("P26").Value = 0.1907052526878442
$ws.Range("Q26").Value = 31.81947539112333
$ws.Range("R26").Value = 286.37527852011
$ws.Range("S26").Value = 0.06853212708821445
$ws.Range("T26").Value = 0.06853212708821445
